$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date format used by column A (matches existing style s="1")
$dateFormat = "MM/DD/YY"
# Time format used by columns B/C (matches existing style s="2"; note the escaped space)
$timeFormat = "HH:MM:SS\ AM/PM"

# --- Row 12 ---
$ws.Range("A12").Value = 41964
$ws.Range("A12").NumberFormat = $dateFormat
$ws.Range("B12").Value = 0.5625
$ws.Range("B12").NumberFormat = $timeFormat
$ws.Range("C12").Value = 0.708333333333333
$ws.Range("C12").NumberFormat = $timeFormat
$ws.Range("D12").Formula = "=ROUND(ABS(C12-B12) * 24, 1)"

# --- Row 13 ---
$ws.Range("A13").Value = 41964
$ws.Range("A13").NumberFormat = $dateFormat
$ws.Range("B13").Value = 0.958333333333333
$ws.Range("B13").NumberFormat = $timeFormat
$ws.Range("C13").Value = 0.999988425925926
$ws.Range("C13").NumberFormat = $timeFormat
$ws.Range("D13").Formula = "=ROUND(ABS(C13-B13) * 24, 1)"

# --- Row 14 ---
$ws.Range("A14").Value = 41965
$ws.Range("A14").NumberFormat = $dateFormat
$ws.Range("B14").Value = 0
$ws.Range("B14").NumberFormat = $timeFormat
$ws.Range("C14").Value = 0.0833333333333333
$ws.Range("C14").NumberFormat = $timeFormat
$ws.Range("D14").Formula = "=ROUND(ABS(C14-B14) * 24, 1)"

# --- Row 15 (new) ---
$ws.Range("A15").Value = 41965
$ws.Range("A15").NumberFormat = $dateFormat
$ws.Range("B15").Value = 0.625
$ws.Range("B15").NumberFormat = $timeFormat
$ws.Range("C15").Value = 0.791666666666667
$ws.Range("C15").NumberFormat = $timeFormat
$ws.Range("D15").Formula = "=ROUND(ABS(C15-B15) * 24, 1)"

# --- Row 16 (new) ---
$ws.Range("A16").Value = 41966
$ws.Range("A16").NumberFormat = $dateFormat
$ws.Range("B16").Value = 0.0416666666666667
$ws.Range("B16").NumberFormat = $timeFormat
$ws.Range("C16").Value = 0.208333333333333
$ws.Range("C16").NumberFormat = $timeFormat
$ws.Range("D16").Formula = "=ROUND(ABS(C16-B16) * 24, 1)"

# --- Row 17 (new) ---
$ws.Range("A17").Value = 41966
$ws.Range("A17").NumberFormat = $dateFormat
$ws.Range("B17").Value = 0.625
$ws.Range("B17").NumberFormat = $timeFormat
$ws.Range("C17").Value = 0.791666666666667
$ws.Range("C17").NumberFormat = $timeFormat
$ws.Range("D17").Formula = "=ROUND(ABS(C17-B17) * 24, 1)"

# --- Total row (formula text tidied up, no leading space) ---
$ws.Range("D26").Formula = "=SUM(D3:D25)"

# --- Update the active selection to D12 ---
$ws.Range("D12").Select()

$wb.Save()
